# Update countries & provincias Spain
# Applies the 2020-05-02 10:03 -> 11:08 data refresh to the "Pais" sheet:
#   - Rusia overtakes Turquia (rows 10/11 swap rank, new Rusia totals)
#   - Eslovaquia overtakes Lituania (rows 85/86 swap rank, new Eslovaquia totals)
#   - refreshed case counts for a handful of other countries
#   - updated "Datos actualizados" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 2 de Mayo de 2020 a las 11:08"

# --- Rank swap: Rusia now above Turquia (row 10 / row 11) -------------
$ws.Range("A10").Value = "Rusia"
$ws.Range("A11").Value = "Turquia"

# --- Rank swap: Eslovaquia now above Lituania (row 85 / row 86) -------
$ws.Range("A85").Value = "Eslovaquia"
$ws.Range("A86").Value = "Lituania"

# --- Refreshed numeric columns (B=Casos totales, C=Nuevos casos, -------
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy,
#     H=Muertes) -----------------------------------------------------
$updates = @{
    "B10" = 124054; "C10" = 9623;  "D10" = 15013; "E10" = 107819; "F10" = 2300; "G10" = 53;  "H10" = 1222;
    "B11" = 122392;                "D11" = 53808; "E11" = 65326;  "F11" = 1480;              "H11" = 3258;

    "B28" = 17548; "C28" = 447;                    "E28" = 16264;

    "B43" = 8928;  "C43" = 156;   "D43" = 1124;  "E43" = 7201;                 "G43" = 24;  "H43" = 603;
    "B44" = 8790;  "C44" = 552;   "D44" = 177;   "E44" = 8438;                 "G44" = 5;   "H44" = 175;
    "B46" = 7740;  "C46" = 3;     "D46" = 3378;  "E46" = 4121;                 "G46" = 1;   "H46" = 241;
    "B49" = 6783;  "C49" = 16;    "D49" = 5789;  "E49" = 901;   "F49" = 28;
    "B51" = 6176;  "C51" = 105;   "D51" = 4326;  "E51" = 1747;  "F51" = 31;
    "B66" = 2483;  "C66" = 36;    "D66" = 750;   "E66" = 1722;

    "D71" = 1229;  "E71" = 856;

    "B85" = 1407;  "C85" = 4;     "D85" = 608;   "E85" = 775;   "F85" = 5;                 "H85" = 24;
    "B86" = 1406;  "C86" = 7;     "D86" = 632;   "E86" = 728;   "F86" = 17;    "G86" = 1;  "H86" = 46;

    "D103" = 172;  "E103" = 511;
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
